$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.535.75'
Set-TextValue 'E2' '  -1.55%  '
Set-TextValue 'D3' '1.880.01'
Set-TextValue 'E3' '  -1.79%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.20%  '
Set-TextValue 'D5' '233.32'
Set-TextValue 'E5' '  -5.21%  '
Set-TextValue 'E6' '  +0.32%  '
Set-TextValue 'D7' '0.4866'
Set-TextValue 'E7' '  -2.01%  '
Set-TextValue 'D8' '0.2876'
Set-TextValue 'E8' '  -3.91%  '
Set-TextValue 'D9' '0.06632'
Set-TextValue 'E9' '  -2.55%  '
Set-TextValue 'D10' '1.878.01'
Set-TextValue 'E10' '  -1.81%  '
Set-TextValue 'D11' '16.66'
Set-TextValue 'E11' '  -2.29%  '
Set-TextValue 'D12' '0.07231'
Set-TextValue 'E12' '  -1.08%  '
Set-TextValue 'D13' '88.17'
Set-TextValue 'E13' '  -2.02%  '
Set-TextValue 'D14' '4.959'
Set-TextValue 'E14' '  -2.28%  '
Set-TextValue 'D15' '0.6599'
Set-TextValue 'E15' '  -3.51%  '
Set-TextValue 'D16' '30.483.68'
Set-TextValue 'E16' '  -1.26%  '
Set-TextValue 'E17' '  +0.10%  '
Set-TextValue 'D18' '0.000007792'
Set-TextValue 'E18' '  -3.00%  '
Set-TextValue 'D19' '12.85'
Set-TextValue 'E19' '  -2.92%  '
Set-TextValue 'D20' '2.120.72'
Set-TextValue 'E20' '  -1.76%  '
Set-TextValue 'E21' '  -0.59%  '
Set-TextValue 'D22' '4.714'
Set-TextValue 'E22' '  -3.46%  '
Set-TextValue 'D23' '185.74'
Set-TextValue 'E23' '  +5.82%  '
Set-TextValue 'D24' '6.009'
Set-TextValue 'E24' '  -1.04%  '
Set-TextValue 'D25' '9.222'
Set-TextValue 'E25' '  -1.22%  '
Set-TextValue 'D26' '155.60'
Set-TextValue 'E26' '  +2.16%  '
Set-TextValue 'D27' '18.32'
Set-TextValue 'E27' '  +1.17%  '
Set-TextValue 'D28' '1.829'
Set-TextValue 'E28' '  -6.22%  '
Set-TextValue 'D29' '1.399'
Set-TextValue 'E29' '  -1.46%  '
Set-TextValue 'D30' '4.224'
Set-TextValue 'E30' '  -2.66%  '
Set-TextValue 'D31' '0.08984'
Set-TextValue 'E31' '  +0.68%  '
Set-TextValue 'D32' '3.901'
Set-TextValue 'E32' '  -4.40%  '
Set-TextValue 'D33' '0.05181'
Set-TextValue 'E33' '  -2.18%  '
Set-TextValue 'D34' '0.7302'
Set-TextValue 'E34' '  -2.63%  '
Set-TextValue 'D35' '1.073'
Set-TextValue 'E35' '  -6.45%  '
Set-TextValue 'D36' '2.700'
Set-TextValue 'E36' '  +1.83%  '
Set-TextValue 'D37' '0.01809'
Set-TextValue 'E37' '  -6.32%  '
Set-TextValue 'D38' '2.644'
Set-TextValue 'E38' '  -3.12%  '
Set-TextValue 'D39' '0.9171'
Set-TextValue 'E39' '  -2.64%  '
Set-TextValue 'D40' '2.027'
Set-TextValue 'E40' '  -8.39%  '
Set-TextValue 'D41' '0.4295'
Set-TextValue 'E41' '  -2.56%  '
Set-TextValue 'D42' '0.9965'
Set-TextValue 'E42' '  -0.44%  '
Set-TextValue 'D43' '103.33'
Set-TextValue 'E43' '  -1.87%  '
Set-TextValue 'D44' '5.666'
Set-TextValue 'E44' '  -5.28%  '
Set-TextValue 'D45' '0.1329'
Set-TextValue 'D46' '7.207'
Set-TextValue 'E46' '  -7.76%  '
Set-TextValue 'D47' '0.05806'
Set-TextValue 'E47' '  -0.69%  '
Set-TextValue 'D48' '8.540'
Set-TextValue 'E48' '  -0.16%  '
Set-TextValue 'D49' '1.398'
Set-TextValue 'E49' '  +1.08%  '
Set-TextValue 'D50' '0.3866'
Set-TextValue 'E50' '  -1.76%  '
Set-TextValue 'D51' '33.05'
Set-TextValue 'E51' '  -1.16%  '
